# edit.ps1 -- apply Daily_Scores update (2025-02-11 recompute + new 2025-02-12 rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 38: I38/P38 values refreshed ---
$ws.Cells.Item(38, 9).Value = 9.853016671285292
$ws.Cells.Item(38, 16).Value = 54.7299622371925

# --- Row 42: recomputed score columns (2025-02-11) ---
$ws.Cells.Item(42, 3).Value = 8.505718531337935
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 9.034194337603026
$ws.Cells.Item(42, 6).Value = 10
$ws.Cells.Item(42, 7).Value = 9.423173432489826
$ws.Cells.Item(42, 8).Value = 10
$ws.Cells.Item(42, 9).Value = 9.807085235056215
$ws.Cells.Item(42, 10).Value = 9.222516160852724
$ws.Cells.Item(42, 11).Value = 7.586106406679217
$ws.Cells.Item(42, 12).Value = 8.480318896884976
$ws.Cells.Item(42, 13).Value = 9.620710346813063
$ws.Cells.Item(42, 14).Value = 0
$ws.Cells.Item(42, 15).Value = 0
$ws.Cells.Item(42, 16).Value = 53.97698828997927
$ws.Cells.Item(42, 17).Value = 37.7028350577377

# --- Row 43: recomputed score columns (2025-02-11) ---
$ws.Cells.Item(43, 3).Value = 6.246646823209494
$ws.Cells.Item(43, 4).Value = 5
$ws.Cells.Item(43, 5).Value = 5.100462125778582
$ws.Cells.Item(43, 6).Value = 8.305148056642892
$ws.Cells.Item(43, 7).Value = 5.601995184038527
$ws.Cells.Item(43, 8).Value = 10
$ws.Cells.Item(43, 9).Value = 6.178075396825397
$ws.Cells.Item(43, 10).Value = 6.25
$ws.Cells.Item(43, 11).Value = 6.344086021505376
$ws.Cells.Item(43, 12).Value = 0
$ws.Cells.Item(43, 13).Value = 6.672650977269684
$ws.Cells.Item(43, 14).Value = 5
$ws.Cells.Item(43, 15).Value = 5
$ws.Cells.Item(43, 16).Value = 41.14391652862706
$ws.Cells.Item(43, 17).Value = 34.55514805664289

# --- New rows 46-49: 2025-02-12 daily scores ---
# Column A holds a date-like label that must stay literal text
# (matches existing A2:A45 cells, which are plain text, not real dates),
# so we prefix with an apostrophe and then strip the resulting style.
$ws.Cells.Item(46, 1).Value = "'2025-02-12"
$ws.Cells.Item(46, 2).Value = "abs_activity"
$ws.Cells.Item(46, 3).Value = 4.905898911705448
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 6.412151258667126
$ws.Cells.Item(46, 6).Value = 9.164496298318326
$ws.Cells.Item(46, 7).Value = 9.218646392694135
$ws.Cells.Item(46, 8).Value = 10
$ws.Cells.Item(46, 9).Value = 9.834180205301827
$ws.Cells.Item(46, 10).Value = 8.726235192354064
$ws.Cells.Item(46, 11).Value = -7.259499591077078
$ws.Cells.Item(46, 12).Value = 10
$ws.Cells.Item(46, 13).Value = 9.082687021086471
$ws.Cells.Item(46, 14).Value = 0
$ws.Cells.Item(46, 15).Value = 0
$ws.Cells.Item(46, 16).Value = 32.19406419837793
$ws.Cells.Item(46, 17).Value = 37.89073149067239

$ws.Cells.Item(47, 1).Value = "'2025-02-12"
$ws.Cells.Item(47, 2).Value = "rel_activity"
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 5
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = 5.279101391449843
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 10
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 9.25
$ws.Cells.Item(47, 11).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 13).Value = 5.454247854233139
$ws.Cells.Item(47, 14).Value = 5
$ws.Cells.Item(47, 15).Value = 5
$ws.Cells.Item(47, 16).Value = 10.45424785423314
$ws.Cells.Item(47, 17).Value = 34.52910139144984

$ws.Cells.Item(48, 1).Value = "'2025-02-12"
$ws.Cells.Item(48, 2).Value = "abs_sleep"
$ws.Cells.Item(48, 3).Value = 8.033333333333333
$ws.Cells.Item(48, 4).Value = 9.133333333333333
$ws.Cells.Item(48, 5).Value = 10
$ws.Cells.Item(48, 6).Value = 5.933333333333334
$ws.Cells.Item(48, 7).Value = 8.699999999999999
$ws.Cells.Item(48, 8).Value = 9.466666666666667
$ws.Cells.Item(48, 9).Value = 10
$ws.Cells.Item(48, 10).Value = 8.333333333333334
$ws.Cells.Item(48, 11).Value = 10
$ws.Cells.Item(48, 12).Value = 10
$ws.Cells.Item(48, 13).Value = 10
$ws.Cells.Item(48, 14).Value = 0
$ws.Cells.Item(48, 15).Value = 0
$ws.Cells.Item(48, 16).Value = 56.73333333333333
$ws.Cells.Item(48, 17).Value = 42.86666666666667

$ws.Cells.Item(49, 1).Value = "'2025-02-12"
$ws.Cells.Item(49, 2).Value = "rel_sleep"
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 10
$ws.Cells.Item(49, 6).Value = 0
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 9.923800700268314
$ws.Cells.Item(49, 9).Value = 9.99232158988257
$ws.Cells.Item(49, 10).Value = 0
$ws.Cells.Item(49, 11).Value = 10
$ws.Cells.Item(49, 12).Value = 10
$ws.Cells.Item(49, 13).Value = 0
$ws.Cells.Item(49, 14).Value = 0
$ws.Cells.Item(49, 15).Value = 0
$ws.Cells.Item(49, 16).Value = 29.99232158988257
$ws.Cells.Item(49, 17).Value = 19.92380070026832

# Strip the auto-applied date-number-format style from column A so the
# new cells stay plain/unstyled like the rest of the date column.
$ws.Range("A46:A49").ClearFormats()

# --- Refresh the sheet dimension to include the new rows ---
$ws.UsedRange | Out-Null
